$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 392, pushing the existing rows 392..471
# down to 393..472, and populate it with the new weekly price record.
$ws.Rows("392").Insert()

$ws.Range("A392").Value = 8
$ws.Range("B392").Value = "Terminal La Palmera de La Serena"
$ws.Range("C392").Value = "Coquimbo"
$ws.Range("D392").Value = 45173
$ws.Range("E392").Value = 4
$ws.Range("F392").Value = 100112012
$ws.Range("G392").Value = "Espinaca"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 1000
$ws.Range("K392").Value = 450
$ws.Range("L392").Value = 500
$ws.Range("M392").Value = 475
$ws.Range("N392").Value = "$/atado 300 a 500 gramos"
$ws.Range("O392").Value = "Provincia del Elquí"
$ws.Range("P392").Value = 950
$ws.Range("Q392").Value = 0.5
$ws.Range("R392").Value = "Hortaliza"
